# Revision of "Loading your Data" documentation sample workbook.
# Update the column headers of the "Untitled.tab" sheet to the new
# Orange-style prefixed names, and move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) text changes.
# Columns C1/D1 (spo-early / spo-mid) keep their original text.
$ws.Range("A1").Value = "mD#function"
$ws.Range("B1").Value = "mS#gene"
$ws.Range("E1").Value = "c#heat 0"
$ws.Range("F1").Value = "i#heat 10"
$ws.Range("G1").Value = "i#heat 20"

# Update the selected / active cell in the sheet view.
$ws.Range("G2").Select() | Out-Null
